# Added the Signup and Signin tests
$wb = $excel.ActiveWorkbook

$signup = $wb.Worksheets.Item(1)
$signup.Name = "signup_data"

# --- fix the typo'd email on the signup sheet -----------------------------
$signup.Range("A2").Value = "user5@example.com"

# Re-create the four hyperlinks so their address/display text matches the
# corrected e-mail. Re-adding always paints the "Hyperlink" cell style, so
# immediately strip that back down to the formatting the source cells
# actually use (plain border+alignment for column A, Excel's default
# Hyperlink look - which is what the template already used - for column B).
$signup.Hyperlinks.Delete()
$signup.Hyperlinks.Add($signup.Range("A2"), "mailto:user5@example.com", "", "mailto:user5@example.com", "user5@example.com")
$signup.Hyperlinks.Add($signup.Range("B2"), "mailto:ValidPass@123", "", "", "ValidPass@123")
$signup.Hyperlinks.Add($signup.Range("A4"), "mailto:another@example.com", "", "mailto:another@example.com", "another@example.com")
$signup.Hyperlinks.Add($signup.Range("B3"), "mailto:Short@123", "", "", "Short@123")

foreach ($addr in @("A2", "A4")) {
    $r = $signup.Range($addr)
    $r.ClearFormats()
    $r.Borders.LineStyle = 1
    $r.HorizontalAlignment = -4108
    $r.VerticalAlignment = -4108
    $r.WrapText = $true
}

# --- new "signin_data" sheet ----------------------------------------------
# Build it as a copy of signup_data so column widths / row heights / cell
# styles / number formats all come along for free, then prune & extend the
# data to match the signin fixture (2 data rows + 2 blank styled rows).
$signup.Copy($null, $signup)
$signin = $wb.Worksheets.Item(2)
$signin.Name = "signin_data"

$signin.Range("A4:C4").ClearContents()
$signin.Hyperlinks.Delete()
$signin.Hyperlinks.Add($signin.Range("A2"), "mailto:user5@example.com", "", "mailto:user5@example.com", "user5@example.com")
$signin.Hyperlinks.Add($signin.Range("B2"), "mailto:ValidPass@123", "", "", "ValidPass@123")
$signin.Hyperlinks.Add($signin.Range("B3"), "mailto:Short@123", "", "", "Short@123")

$r = $signin.Range("A2")
$r.ClearFormats()
$r.Borders.LineStyle = 1
$r.HorizontalAlignment = -4108
$r.VerticalAlignment = -4108
$r.WrapText = $true

$signin.Range("C4").Value = ""
$signin.Range("C4").HorizontalAlignment = -4108
$signin.Range("C4").VerticalAlignment = -4108
$signin.Range("C4").WrapText = $true

$signin.Range("A5:C5").HorizontalAlignment = -4108
$signin.Range("A5:C5").VerticalAlignment = -4108

$signin.Range("A1").ColumnWidth = 34

# --- selection / active sheet ---------------------------------------------
$signup.Range("C7").Select()
$signin.Range("C12").Select()
